# Update the dSF (column F) values for a set of rows on Sheet1.
# These rows were re-pulled / recalculated (mean calculation), so only the
# "final" spread-delta (dSF) column changes; dS0 (column E) stays as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    3  = -4
    8  = -2
    11 = 4
    12 = -8
    14 = -4
    16 = -3
    18 = 5
    21 = 3
    22 = 7
    28 = -2
    32 = 5
    33 = -11
    35 = -3
    40 = -1
    42 = 0
    43 = 1
    45 = -5
    47 = 4
    49 = -1
    50 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
